$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F = "dSF" values being updated (repull data, push all data, mean calculation)
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -2
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 4
